$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header B1: "T/N" -> "SAMPLE_TYPE" with new bold-white-on-orange formatting
$ws.Range("B1").Value = "SAMPLE_TYPE"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Color = 16777215
$ws.Range("B1").Interior.Color = 11851260

# Header F1: "Sample ID" -> "Specimen_Number"
$ws.Range("F1").Value = "Specimen_Number"

# Update selection to F1 (matches post-edit sheet view state)
$ws.Range("F1").Select()
